$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 8.426422666666666
$ws.Cells.Item(2, 8).Value = 25.279268
$ws.Cells.Item(2, 9).Value = 0.1151758588783328
$ws.Cells.Item(2, 10).Value = 0.1151758588783328
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.398034
$ws.Cells.Item(2, 14).Value = 4.194102
$ws.Cells.Item(2, 15).Value = 0.139066772576779
$ws.Cells.Item(2, 16).Value = 0.139066772576779
$ws.Cells.Item(2, 17).Value = 11.78042538637066
$ws.Cells.Item(2, 18).Value = 106.023828477336
$ws.Cells.Item(2, 19).Value = 0.0160171349729683
$ws.Cells.Item(2, 20).Value = 0.0160171349729683

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 8.426422666666666
$ws.Cells.Item(3, 8).Value = 25.279268
$ws.Cells.Item(3, 9).Value = 0.1151758588783328
$ws.Cells.Item(3, 10).Value = 0.1151758588783328
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 7.939250333333333
$ws.Cells.Item(3, 14).Value = 23.817751
$ws.Cells.Item(3, 15).Value = 0.7897418235434783
$ws.Cells.Item(3, 16).Value = 0.7897418235434784
$ws.Cells.Item(3, 17).Value = 66.89947896514087
$ws.Cells.Item(3, 18).Value = 602.0953106862679
$ws.Cells.Item(3, 19).Value = 0.09095919281876086
$ws.Cells.Item(3, 20).Value = 0.09095919281876087

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 8.426422666666666
$ws.Cells.Item(4, 8).Value = 25.279268
$ws.Cells.Item(4, 9).Value = 0.1151758588783328
$ws.Cells.Item(4, 10).Value = 0.1151758588783328
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.715685
$ws.Cells.Item(4, 14).Value = 2.147055
$ws.Cells.Item(4, 15).Value = 0.0711914038797426
$ws.Cells.Item(4, 16).Value = 0.0711914038797426
$ws.Cells.Item(4, 17).Value = 6.030664306193333
$ws.Cells.Item(4, 18).Value = 54.27597875573999
$ws.Cells.Item(4, 19).Value = 0.008199531086603627
$ws.Cells.Item(4, 20).Value = 0.008199531086603627

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 11.78712033333333
$ws.Cells.Item(5, 8).Value = 35.361361
$ws.Cells.Item(5, 9).Value = 0.1611112760180311
$ws.Cells.Item(5, 10).Value = 0.1611112760180311
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.398034
$ws.Cells.Item(5, 14).Value = 4.194102
$ws.Cells.Item(5, 15).Value = 0.139066772576779
$ws.Cells.Item(5, 16).Value = 0.139066772576779
$ws.Cells.Item(5, 17).Value = 16.47879498809133
$ws.Cells.Item(5, 18).Value = 148.309154892822
$ws.Cells.Item(5, 19).Value = 0.0224052251815542
$ws.Cells.Item(5, 20).Value = 0.0224052251815542

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 11.78712033333333
$ws.Cells.Item(6, 8).Value = 35.361361
$ws.Cells.Item(6, 9).Value = 0.1611112760180311
$ws.Cells.Item(6, 10).Value = 0.1611112760180311
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.939250333333333
$ws.Cells.Item(6, 14).Value = 23.817751
$ws.Cells.Item(6, 15).Value = 0.7897418235434783
$ws.Cells.Item(6, 16).Value = 0.7897418235434784
$ws.Cells.Item(6, 17).Value = 93.58089903545678
$ws.Cells.Item(6, 18).Value = 842.228091319111
$ws.Cells.Item(6, 19).Value = 0.1272363129158966
$ws.Cells.Item(6, 20).Value = 0.1272363129158966

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 11.78712033333333
$ws.Cells.Item(7, 8).Value = 35.361361
$ws.Cells.Item(7, 9).Value = 0.1611112760180311
$ws.Cells.Item(7, 10).Value = 0.1611112760180311
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.715685
$ws.Cells.Item(7, 14).Value = 2.147055
$ws.Cells.Item(7, 15).Value = 0.0711914038797426
$ws.Cells.Item(7, 16).Value = 0.0711914038797426
$ws.Cells.Item(7, 17).Value = 8.435865215761668
$ws.Cells.Item(7, 18).Value = 75.92278694185501
$ws.Cells.Item(7, 19).Value = 0.01146973792058034
$ws.Cells.Item(7, 20).Value = 0.01146973792058034

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 52.94781866666667
$ws.Cells.Item(8, 8).Value = 158.843456
$ws.Cells.Item(8, 9).Value = 0.7237128651036362
$ws.Cells.Item(8, 10).Value = 0.7237128651036362
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.398034
$ws.Cells.Item(8, 14).Value = 4.194102
$ws.Cells.Item(8, 15).Value = 0.139066772576779
$ws.Cells.Item(8, 16).Value = 0.139066772576779
$ws.Cells.Item(8, 17).Value = 74.02285072183467
$ws.Cells.Item(8, 18).Value = 666.205656496512
$ws.Cells.Item(8, 19).Value = 0.1006444124222565
$ws.Cells.Item(8, 20).Value = 0.1006444124222565

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 52.94781866666667
$ws.Cells.Item(9, 8).Value = 158.843456
$ws.Cells.Item(9, 9).Value = 0.7237128651036362
$ws.Cells.Item(9, 10).Value = 0.7237128651036362
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 7.939250333333333
$ws.Cells.Item(9, 14).Value = 23.817751
$ws.Cells.Item(9, 15).Value = 0.7897418235434783
$ws.Cells.Item(9, 16).Value = 0.7897418235434784
$ws.Cells.Item(9, 17).Value = 420.3659869986062
$ws.Cells.Item(9, 18).Value = 3783.293882987456
$ws.Cells.Item(9, 19).Value = 0.571546317808821
$ws.Cells.Item(9, 20).Value = 0.571546317808821

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 52.94781866666667
$ws.Cells.Item(10, 8).Value = 158.843456
$ws.Cells.Item(10, 9).Value = 0.7237128651036362
$ws.Cells.Item(10, 10).Value = 0.7237128651036362
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.715685
$ws.Cells.Item(10, 14).Value = 2.147055
$ws.Cells.Item(10, 15).Value = 0.0711914038797426
$ws.Cells.Item(10, 16).Value = 0.0711914038797426
$ws.Cells.Item(10, 17).Value = 37.89395960245334
$ws.Cells.Item(10, 18).Value = 341.04563642208
$ws.Cells.Item(10, 19).Value = 0.05152213487255864
$ws.Cells.Item(10, 20).Value = 0.05152213487255864
